$wb = $excel.ActiveWorkbook

$wsDemo3 = $wb.Worksheets.Item("DEMO3")
$wsDemo5 = $wb.Worksheets.Item("DEMO5")

$wsDemo3.Range("A5:BU8").Copy($wsDemo5.Range("A11"))
$wsDemo3.Range("A5:BU8").Clear()
$wsDemo3.Rows("5:8").RowHeight = 15
